$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Physiology")

# Insert a brand-new row at position 12 (shifts old rows 12..101 down to 13..102),
# matching the formatting of row 11 (the row immediately above), since the new
# record logically continues the "2016_07_21_0004" entry directly above it
# (noted as "Can be concatenated." in that row).
$ws.Rows.Item(12).Insert()

# Copy the full formatting (borders/alignment/font/fill/number format) from row 11
# onto the freshly inserted row 12.
$ws.Range("A11:I11").Copy()
$ws.Range("A12:I12").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

# Row height and the G:I merge aren't carried by PasteSpecial, so set them explicitly.
$ws.Rows.Item(12).RowHeight = $ws.Rows.Item(11).RowHeight
$ws.Range("G12:I12").Merge()

# Populate the new row's values.
$ws.Cells.Item(12, 1).Value = 1
$ws.Cells.Item(12, 2).Value = "2016_07_21_concat3and4"
$ws.Cells.Item(12, 3).Value = "V"
$ws.Cells.Item(12, 4).Value = 5
$ws.Cells.Item(12, 5).Value = 34
$ws.Cells.Item(12, 6).Value = -87
$ws.Cells.Item(12, 7).Value = "concatenated file for analysis"

# Reflect the scroll position / selection the author ended up with after the edit.
$ws.Application.ActiveWindow.ScrollRow = 13
$ws.Range("G15:I15").Select()
